# Actualización automática hashcode vie feb 22 01:42:10 CET 2019
# Update the MD5-looking hashcode values (column B) for the rows whose
# key (column A) identifies them, as described in the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B100").Value = "aed0b10bec2e9f418199ec1ba1e7d195"
$ws.Range("B104").Value = "d3250a5c1d6eed791df60eeb4e0dcd1e"
$ws.Range("B113").Value = "d29b6e376f4ab6820f3ee6102e491c52"
$ws.Range("B122").Value = "ee278e6bca7a8de6505b4498ce294b23"
$ws.Range("B164").Value = "2c5adae7a14dbf122a3e7e333946951c"
$ws.Range("B230").Value = "6cc38a03a89a547d65027e64ed10d11e"
$ws.Range("B233").Value = "bf8a8bb894e8abfbc617dca6f524925d"
$ws.Range("B331").Value = "e40c86b9f34ec9b35c007636ab13d92c"
$ws.Range("B342").Value = "b105fc2bcbf5ba38e657e44c9d60c42d"
$ws.Range("B343").Value = "9683867abb9ed08c897898b9ce16f688"
$ws.Range("B419").Value = "afba4ee92bb44bede48ddf483ac24705"
$ws.Range("B619").Value = "e3ee95ef384d09352f2806899d18ac19"
$ws.Range("B623").Value = "ff9f888e91bca8d85efafc7661513a32"
$ws.Range("B628").Value = "846c9647ded4ae397a5a92e7ec1d0301"
$ws.Range("B757").Value = "4b38c4c301a7b483c17e2da034918036"
$ws.Range("B760").Value = "fd27f7cf29ffa4a7b917c79611b0d651"
$ws.Range("B763").Value = "e66db9edeb85723f367334b05a32f91c"
$ws.Range("B767").Value = "e2db6fe7b898e2dc1173466e6c23a763"
$ws.Range("B779").Value = "d7f4356c35eb2b8b0deaac7d4e0be00c"
$ws.Range("B818").Value = "1dcbd17e31672161575ab11d3dad0626"
$ws.Range("B831").Value = "572bd04638f6b5cead7fe5e2de230d72"
